$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cell values are stored as plain text (coinranking price/volume strings).
# Column D values that parse as plain numbers need a leading apostrophe so
# Excel keeps them as text instead of silently converting them to numbers
# (this matches how the source data keeps thousands-dot-separated / exotic
# values as text already). Column E values (padded "  +x.xx%  " strings) are
# never auto-numeric so they do not need the prefix.

$ws.Range("D2").Value = '34.555.22'
$ws.Range("E2").Value = '  +2.37%  '
$ws.Range("D3").Value = '1.787.04'
$ws.Range("E3").Value = '  +0.84%  '
$ws.Range("E4").Value = '  +0.11%  '
$ws.Range("D5").Value = '''223.68'
$ws.Range("E5").Value = '  -0.45%  '
$ws.Range("E6").Value = '  +1.00%  '
$ws.Range("E7").Value = '  +0.16%  '
$ws.Range("D8").Value = '''32.97'
$ws.Range("E8").Value = '  +7.99%  '
$ws.Range("E9").Value = '  +1.14%  '
$ws.Range("D10").Value = '''0.0679'
$ws.Range("E10").Value = '  +2.92%  '
$ws.Range("D11").Value = '''0.0935'
$ws.Range("E11").Value = '  +1.47%  '
$ws.Range("D12").Value = '2.043.11'
$ws.Range("E12").Value = '  +0.91%  '
$ws.Range("D13").Value = '''11.08'
$ws.Range("E13").Value = '  +11.02%  '
$ws.Range("D14").Value = '1.781.19'
$ws.Range("E14").Value = '  +0.47%  '
$ws.Range("E15").Value = '  +0.94%  '
$ws.Range("D16").Value = '34.541.36'
$ws.Range("E16").Value = '  +2.43%  '
$ws.Range("E17").Value = '  +2.77%  '
$ws.Range("D18").Value = '''68.53'
$ws.Range("E18").Value = '  +0.28%  '
$ws.Range("D19").Value = '''253.49'
$ws.Range("E19").Value = '  +0.88%  '
$ws.Range("D20").Value = '0.0₃0777'
$ws.Range("E20").Value = '  +5.49%  '
$ws.Range("E21").Value = '  -0.02%  '
$ws.Range("D22").Value = '''10.43'
$ws.Range("E22").Value = '  +1.81%  '
$ws.Range("E23").Value = '  +1.32%  '
$ws.Range("E24").Value = '  +0.31%  '
$ws.Range("D25").Value = '''159.14'
$ws.Range("E25").Value = '  +0.62%  '
$ws.Range("D26").Value = '''16.36'
$ws.Range("E26").Value = '  -0.61%  '
$ws.Range("E27").Value = '  +2.54%  '
$ws.Range("E28").Value = '  +0.10%  '
$ws.Range("E29").Value = '  +0.12%  '
$ws.Range("E30").Value = '  -0.91%  '
$ws.Range("D31").Value = '''0.0515'
$ws.Range("E31").Value = '  +0.72%  '
$ws.Range("E32").Value = '  +0.08%  '
$ws.Range("E33").Value = '  +1.27%  '
$ws.Range("E34").Value = '  +2.92%  '
$ws.Range("D35").Value = '1.445.47'
$ws.Range("E36").Value = '  -0.57%  '
$ws.Range("E37").Value = '  +2.49%  '
$ws.Range("E38").Value = '  -0.82%  '
$ws.Range("D39").Value = '''83.17'
$ws.Range("E39").Value = '  -0.05%  '
$ws.Range("E40").Value = '  +4.77%  '
$ws.Range("E41").Value = '  +0.14%  '
$ws.Range("D42").Value = '''0.897'
$ws.Range("E42").Value = '  +1.34%  '
$ws.Range("D43").Value = '''2.07'
$ws.Range("E43").Value = '  -0.40%  '
$ws.Range("D44").Value = '''0.0502'
$ws.Range("E44").Value = '  -2.04%  '
$ws.Range("D45").Value = '''5.90'
$ws.Range("E46").Value = '  -2.49%  '
$ws.Range("D47").Value = '1.940.68'
$ws.Range("E47").Value = '  +0.98%  '
$ws.Range("D48").Value = '''104.34'
$ws.Range("E48").Value = '  +7.07%  '
$ws.Range("E49").Value = '  +0.08%  '
$ws.Range("D50").Value = '''11.94'
$ws.Range("E50").Value = '  +1.33%  '
$ws.Range("D51").Value = '''49.30'
$ws.Range("E51").Value = '  -2.47%  '
